$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: username becomes "adm" (was "Admin")
$ws.Range("A2").Value = "adm"

# Row 3: username becomes "Admin" (was "adm"); password becomes text "admin123" (was numeric 123)
$ws.Range("A3").Value = "Admin"
$ws.Range("B3").ClearFormats()
$ws.Range("B3").Value = "admin123"

# Row 4: replace with what used to be row 5's contents
$ws.Range("A4").Value = "qasmart.zee@bssuniversal.com"
$ws.Range("B4").Value = "Bss@2025-1"

# Delete old row 5 (its content has been moved up into row 4)
$ws.Rows("5").Delete()

$ws.Range("A2").Select()
